# Lesson 2 slides: minor typo fix — add the missing trailing period to
# "Use raw_input() to allow a user to type a DC address" on slide 20.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(20)
$shape = $s.Shapes.Item(1)
$tr = $shape.TextFrame.TextRange

# 4th paragraph in the body placeholder is the "Use raw_input() ... DC address" line.
$para4 = $tr.Paragraphs(4, 1)
$newRun = $para4.InsertAfter(".")
